# add image in attendance sheet
# Adds an "image" column (E) to the attendance sheet: each attendee gets a
# profile-photo URL, three of which are turned into real hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Candidate profile-photo URLs (freepik stock photos)
$img1 = "https://img.freepik.com/free-photo/handsome-young-man-with-arms-crossed-white-background_23-2148222620.jpg"
$img2 = "https://img.freepik.com/free-photo/handsome-confident-smiling-man-with-hands-crossed-chest_176420-18743.jpg"
$img3 = "https://img.freepik.com/free-photo/attractive-mixed-race-male-with-positive-smile-shows-white-teeth-keeps-hands-stomach-being-high-spirit-wears-white-shirt-rejoices-positive-moments-life-people-emotions-concept_273609-15527.jpg"
$img4 = "https://img.freepik.com/free-photo/confident-handsome-guy-posing-against-white-wall_176420-32936.jpg"
$img5 = "https://img.freepik.com/free-photo/fashionable-stylish-man-with-dark-eyes-casual-clothes-looking-aside-with-placid-thoughtful-look-pensive-guy-with-puzzled-expression-thinking-about-something-building-plans_176420-10331.jpg"
$img6 = "https://img.freepik.com/free-photo/thoughtful-concerned-man-thinking-trying-find-solution_176420-19574.jpg"
$img7 = "https://img.freepik.com/free-photo/serious-thoughtful-man-squinting-skeptical-thinking-as-making-choice_176420-19020.jpg"

# Header
$ws.Range("E1").Value = "image"

# Fill the new column in row order, so the shared-string table is built up
# in the same order Excel itself would encounter the values.
$ws.Range("E2").Value = $img1
$ws.Range("E3").Value = $img2
$ws.Range("E4").Value = $img3
$ws.Range("E5").Value = $img4
$ws.Range("E6").Value = $img5
$ws.Range("E7").Value = $img6
$ws.Range("E8").Value = $img7
$ws.Range("E9").Value = $img3
$ws.Range("E10").Value = $img6
$ws.Range("E11").Value = $img1

# Turn three of the cells into real (clickable) hyperlinks.
$ws.Hyperlinks.Add($ws.Range("E4"), $img3)
$ws.Hyperlinks.Add($ws.Range("E7"), $img6)
$ws.Hyperlinks.Add($ws.Range("E2"), $img1)

# Match the selection left behind in the saved workbook.
$ws.Range("E1:E11").Select()
